$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear rows that no longer hold content in the final layout
$ws.Range("B43").ClearContents()
$ws.Range("B45").ClearContents()
$ws.Range("B51").ClearContents()

# Write final text content per row (column B)
$ws.Range("B42").Value = 'Khi tự động hóa trình duyệt Internet Explorer (IE) trong môi trường Windows Server 2019 Datacenter, có những ưu điểm và nhược điểm cần xem xét:'
$ws.Range("B44").Value = 'Ưu điểm:'
$ws.Range("B46").Value = 'Quản lý tập trung: Môi trường Windows Server 2019 Datacenter thường được sử dụng trong các hệ thống doanh nghiệp lớn, có khả năng quản lý tập trung cao hơn cho việc triển khai và duy trì các kịch bản tự động hóa.'
$ws.Range("B47").Value = 'Hỗ trợ cao cấp: Windows Server 2019 Datacenter cung cấp nhiều tính năng và dịch vụ cao cấp, giúp cho việc triển khai và quản lý các ứng dụng tự động hóa trở nên hiệu quả hơn.'
$ws.Range("B48").Value = 'Hiệu suất: Môi trường máy chủ thường có khả năng xử lý và hiệu suất cao hơn so với máy tính cá nhân, giúp tự động hóa trình duyệt IE chạy nhanh hơn và ổn định hơn.'
$ws.Range("B49").Value = 'Phân quyền và bảo mật: Môi trường máy chủ cung cấp khả năng phân quyền và quản lý bảo mật nâng cao, giúp đảm bảo rằng quyền truy cập và thực thi của mã tự động hóa được kiểm soát cẩn thận.'
$ws.Range("B50").Value = 'Tích hợp hệ thống: Tự động hóa trình duyệt IE trong môi trường máy chủ có thể dễ dàng tích hợp với các tác vụ hệ thống khác và quản lý toàn bộ quá trình một cách tốt hơn.'
$ws.Range("B52").Value = 'Nhược điểm:'
$ws.Range("B54").Value = 'Khả năng tương thích: IE không còn là trình duyệt phát triển chính thống, và không hỗ trợ nhiều tính năng và tiêu chuẩn web hiện đại. Việc tự động hóa trình duyệt IE trong môi trường Windows Server 2019 Datacenter có thể gặp khó khăn trong việc đảm bảo tương thích.'
$ws.Range("B55").Value = 'Hiệu suất không đảm bảo: Mặc dù máy chủ có hiệu suất tốt hơn, việc tự động hóa trình duyệt IE vẫn có thể ảnh hưởng đến hiệu suất toàn hệ thống. Việc kiểm soát hiệu suất là điều cần quan tâm.'
$ws.Range("B56").Value = 'Phức tạp hóa quản lý: Trong môi trường máy chủ, việc triển khai và quản lý các kịch bản tự động hóa có thể phức tạp hơn do tính chất phức tạp của hệ thống.'
$ws.Range("B57").Value = 'Khả năng gỡ lỗi: Môi trường máy chủ có thể không cung cấp các công cụ gỡ lỗi và môi trường phát triển tích hợp mạnh mẽ như máy tính cá nhân.'
$ws.Range("B59").Value = 'Tóm lại, tự động hóa trình duyệt IE trong môi trường Windows Server 2019 Datacenter có thể mang lại lợi ích về quản lý, hiệu suất và tích hợp hệ thống. Tuy nhiên, cần xem xét cẩn thận về tương thích và bảo mật, và xem xét việc sử dụng các công cụ tự động hóa hiện đại hơn để đảm bảo tính hiệu quả và bảo mật của quá trình tự động hóa.'

# Update view state: scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("E45").Select() | Out-Null
